$d = $word.ActiveDocument

# 1) Status cell text change: "Por hacer" -> "Completado"
$d.Content.Find.Execute("Por hacer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Completado", 2) | Out-Null

# 2) Table #3 (1-based Tables.Item(3)) column width tweaks:
#    gridCol 538 -> 537 (col 1) and 1336 -> 1337 (col 4), applied per-row
#    since this table has no merged cells, per-cell Width assignment
#    updates both the tblGrid and every row's tcW correctly.
$t3 = $d.Tables.Item(3)
for ($r = 1; $r -le $t3.Rows.Count; $r++) {
    $t3.Cell($r, 1).Width = 26.85   # 537 dxa
    $t3.Cell($r, 4).Width = 66.85   # 1337 dxa
}

# 3) Table #6 (1-based Tables.Item(6)) column width tweaks:
#    gridCol 3150 -> 3149 (col 2) and 1892 -> 1893 (col 3), applied per-row
#    for the regular (non-merged) rows 1-3. Row 4 has a merged cell
#    (gridSpan=2) covering columns 1-2, so its last cell shares the
#    engine's internal column-2 slot with rows 1-3's column 2 cells;
#    column 2 is written last so the shared slot ends up holding the
#    value rows 1-3 need.
$t6 = $d.Tables.Item(6)
for ($r = 1; $r -le 3; $r++) {
    $t6.Cell($r, 3).Width = 94.65   # 1893 dxa
}
for ($r = 1; $r -le 3; $r++) {
    $t6.Cell($r, 2).Width = 157.45  # 3149 dxa
}
